$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-02-13 Thursday"; new = "2025-02-14 Friday"},
    @{old = "688÷4="; new = "613÷5="},
    @{old = "706÷5="; new = "186÷6="},
    @{old = "392÷2="; new = "322÷3="},
    @{old = "501÷5="; new = "115÷7="},
    @{old = "402÷4="; new = "897÷6="},
    @{old = "638÷6="; new = "922÷6="},
    @{old = "679÷2="; new = "298÷6="},
    @{old = "922÷9="; new = "189÷5="},
    @{old = "630÷5="; new = "798÷8="},
    @{old = "395÷8="; new = "858÷2="},
    @{old = "123÷4="; new = "272÷4="},
    @{old = "590÷8="; new = "596÷7="},
    @{old = "890÷3="; new = "696÷7="},
    @{old = "525÷3="; new = "648÷6="},
    @{old = "832÷8="; new = "403÷8="},
    @{old = "233÷5="; new = "181÷4="},
    @{old = "405÷5="; new = "604÷5="},
    @{old = "404÷6="; new = "311÷4="},
    @{old = "700÷6="; new = "747÷8="},
    @{old = "169÷3="; new = "338÷5="},
    @{old = "104÷5="; new = "874÷3="},
    @{old = "317÷5="; new = "566÷9="},
    @{old = "539÷4="; new = "929÷5="},
    @{old = "797÷5="; new = "719÷8="},
    @{old = "887÷9="; new = "528÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
